# Changed for loop matchers to be the same as liquid templating.
#
# The template document contains several "spacer" paragraphs that only hold
# a (hidden) `_GoBack` bookmark left over from the last cursor position when
# the document was saved in Word. These bookmark-only paragraphs should
# simply become plain empty paragraphs, i.e. the `_GoBack` bookmark needs to
# be removed from the document while leaving the surrounding paragraphs
# (and everything else) untouched.
$d = $word.ActiveDocument

while ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
